$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.856.45"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.808.94"
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.15"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4466"
$ws.Range("E7").Value = "  +5.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3669"
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07276"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.63"
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.817.67"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.596"
$ws.Range("E13").Value = "  -2.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07065"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.291"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.41"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008712"
$ws.Range("E18").Value = "  -1.97%  "
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.875.15"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.138"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  -1.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.980"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.18"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.219"
$ws.Range("E26").Value = "  +1.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.33"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.191"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.91"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08810"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("E31").Value = "  -3.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7454"
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.929"
$ws.Range("E33").Value = "  +3.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.430"
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.000"
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.089"
$ws.Range("E36").Value = "  -3.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01956"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05174"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5271"
$ws.Range("E39").Value = "  +3.58%  "
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.048"
$ws.Range("E41").Value = "  -3.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1685"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5158"
$ws.Range("E43").Value = "  +7.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.422"
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.970"
$ws.Range("E45").Value = "  +6.84%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.47"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.31"
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9997"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.649"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06319"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9130"
$ws.Range("E51").Value = "  -0.71%  "
